$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checkertificate")

# Row 12
$ws.Range("A12").Value = "01.07.2023 01:45 (Kyiv+Israel) 23:45 (UTC) 08:45 (Japan) 05:15 (India)"
$ws.Range("B12").Value = 0.616
$ws.Range("C12").Value = 0.131
$ws.Range("D12").Value = "***"
$ws.Range("E12").Value = "***"

# Row 13
$ws.Range("A13").Value = "01.07.2023 01:47 (Kyiv+Israel) 23:47 (UTC) 08:47 (Japan) 05:17 (India)"
$ws.Range("B13").Value = "***"
$ws.Range("C13").Value = "***"
$ws.Range("D13").Value = 0.6840000000000001
$ws.Range("E13").Value = 0.201
